$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 119.5
$ws.Range("I4").Value = 119.5
$ws.Range("K4").Value = 119.5
$ws.Range("M4").Value = -5.5

$ws.Range("H17").Value = 3112.0588
$ws.Range("J17").Value = 3112.0588
$ws.Range("L17").Value = 9336.1764
$ws.Range("N17").Value = -9672.1764

$ws.Range("H92").Value = 613.46155
$ws.Range("I92").Value = 636.9167
$ws.Range("K92").Value = 636.9167
$ws.Range("M92").Value = 611.0833

$ws.Range("H132").Value = 7836.81
$ws.Range("I132").Value = 1755.6
$ws.Range("J132").Value = 12812.346
$ws.Range("K132").Value = 5266.799999999999
$ws.Range("L132").Value = 38437.038
$ws.Range("M132").Value = -2736.799999999999
$ws.Range("N132").Value = -43497.038

$ws.Range("H138").Value = 5488.082
$ws.Range("I138").Value = 2357.7036
$ws.Range("J138").Value = 7973.9707
$ws.Range("K138").Value = 7073.110799999999
$ws.Range("L138").Value = 23921.9121
$ws.Range("M138").Value = -1933.110799999999
$ws.Range("N138").Value = -34201.9121

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6156.6924
$ws.Range("I32").Value = 3827.75
$ws.Range("K32").Value = 3827.75
$ws.Range("M32").Value = -3540.75

$ws.Range("H45").Value = 1934.7778
$ws.Range("I45").Value = 1703
$ws.Range("J45").Value = 2398.3333
$ws.Range("K45").Value = 1703
$ws.Range("L45").Value = 2398.3333
$ws.Range("M45").Value = -1326
$ws.Range("N45").Value = -3152.3333

$ws.Range("H132").Value = 15667.782
$ws.Range("I132").Value = 19574.709
$ws.Range("K132").Value = 58724.12699999999
$ws.Range("M132").Value = -56194.12699999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2172.6155
$ws.Range("I20").Value = 2317.2222
$ws.Range("J20").Value = 1847.25
$ws.Range("K20").Value = 2317.2222
$ws.Range("L20").Value = 1847.25
$ws.Range("M20").Value = -2070.2222
$ws.Range("N20").Value = -2341.25

$ws.Range("H99").Value = 1159421.5
$ws.Range("I99").Value = 1303974.2
$ws.Range("K99").Value = 1303974.2
$ws.Range("M99").Value = -1302476.2

$ws.Range("H129").Value = 52500
$ws.Range("J129").Value = 52500
$ws.Range("L129").Value = 52500
$ws.Range("N129").Value = -62500

$ws.Range("H130").Value = 62709
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null

$ws.Range("H134").Value = 2189.7368
$ws.Range("I134").Value = 2052.7334
$ws.Range("J134").Value = 2703.5
$ws.Range("K134").Value = 6158.2002
$ws.Range("L134").Value = 8110.5
$ws.Range("M134").Value = -3623.2002
$ws.Range("N134").Value = -13180.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3758.4285
$ws.Range("I31").Value = 3161.6511
$ws.Range("J31").Value = 5732.385
$ws.Range("K31").Value = 3161.6511
$ws.Range("L31").Value = 5732.385
$ws.Range("M31").Value = -2866.6511
$ws.Range("N31").Value = -6322.385

$ws.Range("H34").Value = 3758.4285
$ws.Range("I34").Value = 3161.6511
$ws.Range("J34").Value = 5732.385
$ws.Range("K34").Value = 3161.6511
$ws.Range("L34").Value = 5732.385
$ws.Range("M34").Value = -2959.6511
$ws.Range("N34").Value = -6136.385

$ws.Range("H58").Value = 436267.88
$ws.Range("I58").Value = 626310.6
$ws.Range("K58").Value = 626310.6
$ws.Range("M58").Value = -626107.6

$ws.Range("H136").Value = 436267.88
$ws.Range("I136").Value = 626310.6
$ws.Range("K136").Value = 1878931.8
$ws.Range("M136").Value = -1876381.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 36507584
$ws.Range("I4").Value = 1083647.8
$ws.Range("K4").Value = 3250943.4
$ws.Range("M4").Value = -3250831.4

$ws.Range("H5").Value = 1214
$ws.Range("I5").Value = 1016.8
$ws.Range("J5").Value = 2200
$ws.Range("K5").Value = 3050.4
$ws.Range("L5").Value = 6600
$ws.Range("M5").Value = -2938.4
$ws.Range("N5").Value = -6824

$ws.Range("H38").Value = 166.1
$ws.Range("J38").Value = 323.375
$ws.Range("L38").Value = 970.125
$ws.Range("N38").Value = -1664.125

$ws.Range("H97").Value = 291
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = $null

$ws.Range("H131").Value = 6946383.5
$ws.Range("J131").Value = 4977153
$ws.Range("L131").Value = 14931459
$ws.Range("N131").Value = -14941539

$ws.Range("H135").Value = 1214
$ws.Range("I135").Value = 1016.8
$ws.Range("J135").Value = 2200
$ws.Range("K135").Value = 9151.199999999999
$ws.Range("L135").Value = 19800
$ws.Range("M135").Value = -6616.199999999999
$ws.Range("N135").Value = -24870

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3739
$ws.Range("I126").Value = 2154.4443
$ws.Range("J126").Value = 6908.1113
$ws.Range("K126").Value = 6463.3329
$ws.Range("L126").Value = 20724.3339
$ws.Range("M126").Value = -3993.3329
$ws.Range("N126").Value = -25664.3339

$ws.Range("H132").Value = 6662.59
$ws.Range("I132").Value = 6074.5483
$ws.Range("K132").Value = 18223.6449
$ws.Range("M132").Value = -15693.6449

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = $null

$ws.Range("H40").Value = 3511.2942
$ws.Range("I40").Value = 2958.6667
$ws.Range("K40").Value = 2958.6667
$ws.Range("M40").Value = -2822.6667

$ws.Range("H68").Value = 949186.9399999999
$ws.Range("I68").Value = 1749975.9
$ws.Range("J68").Value = 2799.9092
$ws.Range("K68").Value = 1749975.9
$ws.Range("L68").Value = 2799.9092
$ws.Range("M68").Value = -1749226.9
$ws.Range("N68").Value = -4297.9092

$ws.Range("H71").Value = 949186.9399999999
$ws.Range("I71").Value = 1749975.9
$ws.Range("J71").Value = 2799.9092
$ws.Range("K71").Value = 8749879.5
$ws.Range("L71").Value = 13999.546
$ws.Range("M71").Value = -8746135.5
$ws.Range("N71").Value = -21487.546

$ws.Range("H132").Value = 6661.385
$ws.Range("I132").Value = 4875.4443
$ws.Range("J132").Value = 7606.8823
$ws.Range("K132").Value = 14626.3329
$ws.Range("L132").Value = 22820.6469
$ws.Range("M132").Value = -12096.3329
$ws.Range("N132").Value = -27880.6469

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5036.909
$ws.Range("I62").Value = 3822.889
$ws.Range("J62").Value = 10500
$ws.Range("K62").Value = 3822.889
$ws.Range("L62").Value = 10500
$ws.Range("M62").Value = -3198.889
$ws.Range("N62").Value = -11748

$ws.Range("H65").Value = 5036.909
$ws.Range("I65").Value = 3822.889
$ws.Range("J65").Value = 10500
$ws.Range("K65").Value = 19114.445
$ws.Range("L65").Value = 52500
$ws.Range("M65").Value = -15994.445
$ws.Range("N65").Value = -58740

$ws.Range("H122").Value = 2010.7764
$ws.Range("I122").Value = 1968.381
$ws.Range("J122").Value = 2216.2307
$ws.Range("K122").Value = 5905.143
$ws.Range("L122").Value = 6648.6921
$ws.Range("M122").Value = -3455.143
$ws.Range("N122").Value = -11548.6921

$ws.Range("H126").Value = 1265.1666
$ws.Range("I126").Value = 1220.7778
$ws.Range("K126").Value = 3662.3334
$ws.Range("M126").Value = -1192.3334

$ws.Range("H132").Value = 13894488
$ws.Range("I132").Value = 1690.2222
$ws.Range("J132").Value = 55572884
$ws.Range("K132").Value = 5070.6666
$ws.Range("L132").Value = 166718652
$ws.Range("M132").Value = -2540.6666
$ws.Range("N132").Value = -166723712

$ws.Range("H136").Value = 7016.37
$ws.Range("I136").Value = 1893.5172
$ws.Range("J136").Value = 9108.803
$ws.Range("K136").Value = 5680.5516
$ws.Range("L136").Value = 27326.409
$ws.Range("M136").Value = -3130.5516
$ws.Range("N136").Value = -32426.409
